$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking strings (e.g. "36.538.12", "0.380")
# that must remain literal text (matching the source inlineStr cells), so
# force text format before assignment to stop Excel auto-converting them to
# floating point numbers (which would also eat significant trailing zeros).
$textCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D30", "D32", "D33", "D34", "D35", "D39", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "36.538.12"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "2.100.89"
$ws.Range("E3").Value = "  +9.67%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "253.41"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").Value = "0.657"
$ws.Range("E6").Value = "  -6.67%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "47.71"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").Value = "0.380"
$ws.Range("E9").Value = "  +2.16%  "
$ws.Range("D10").Value = "59.93"
$ws.Range("E10").Value = "  +2.37%  "
$ws.Range("D11").Value = "0.0743"
$ws.Range("E11").Value = "  -2.83%  "
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "14.58"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").Value = "2.406.10"
$ws.Range("E14").Value = "  +9.67%  "
$ws.Range("D15").Value = "0.838"
$ws.Range("E15").Value = "  +2.83%  "
$ws.Range("D16").Value = "2.099.34"
$ws.Range("E16").Value = "  +9.48%  "
$ws.Range("D17").Value = "5.10"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "36.559.18"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "72.83"
$ws.Range("E19").Value = "  -2.89%  "
$ws.Range("E20").Value = "  -3.70%  "
$ws.Range("D21").Value = "13.26"
$ws.Range("E21").Value = "  -1.21%  "
$ws.Range("D22").Value = "240.38"
$ws.Range("E22").Value = "  -4.32%  "
$ws.Range("D23").Value = "5.24"
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "2.53"
$ws.Range("E25").Value = "  -4.02%  "
$ws.Range("D26").Value = "171.46"
$ws.Range("E26").Value = "  +1.95%  "
$ws.Range("D27").Value = "21.38"
$ws.Range("E27").Value = "  +13.85%  "
$ws.Range("D28").Value = "9.17"
$ws.Range("E28").Value = "  +4.06%  "
$ws.Range("E29").Value = "  -9.42%  "
$ws.Range("D30").Value = "28.84"
$ws.Range("E30").Value = "  +60.88%  "
$ws.Range("E31").Value = "  -5.32%  "
$ws.Range("D32").Value = "4.50"
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("D33").Value = "0.0618"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "2.46"
$ws.Range("E34").Value = "  +21.95%  "
$ws.Range("D35").Value = "0.987"
$ws.Range("E35").Value = "  +12.62%  "
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("D39").Value = "4.09"
$ws.Range("E39").Value = "  -5.98%  "
$ws.Range("E40").Value = "  -11.18%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "1.17"
$ws.Range("E41").Value = "  +6.30%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0224"
$ws.Range("E42").Value = "  -1.92%  "
$ws.Range("D43").Value = "98.23"
$ws.Range("E43").Value = "  -6.78%  "
$ws.Range("D45").Value = "16.02"
$ws.Range("E45").Value = "  -9.51%  "
$ws.Range("D46").Value = "1.333.02"
$ws.Range("E46").Value = "  -1.30%  "
$ws.Range("D47").Value = "0.0842"
$ws.Range("E47").Value = "  +3.59%  "
$ws.Range("D48").Value = "7.09"
$ws.Range("E48").Value = "  +9.62%  "
$ws.Range("D49").Value = "2.306.59"
$ws.Range("E49").Value = "  +10.26%  "
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("E51").Value = "  -5.65%  "
